$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 173 (shifts old rows 173-252 down to 174-253)
$ws.Range("A173").EntireRow.Insert()

# Populate the new row 173 with its full record (same as the surrounding
# weekly records for this market/product, with the new week's data)
$ws.Range("A173").Value2 = 11
$ws.Range("B173").Value2 = "Vega Monumental Concepción"
$ws.Range("C173").Value2 = "Bíobío"
$ws.Range("D173").Value2 = 45134
$ws.Range("E173").Value2 = 8
$ws.Range("F173").Value2 = 100112032
$ws.Range("G173").Value2 = "Zapallo italiano"
$ws.Range("H173").Value2 = "Sin especificar"
$ws.Range("I173").Value2 = "Primera"
$ws.Range("J173").Value2 = 100
$ws.Range("K173").Value2 = 15000
$ws.Range("L173").Value2 = 15000
$ws.Range("M173").Value2 = 15000
$ws.Range("N173").Value2 = "$/caja 50 unidades"
$ws.Range("O173").Value2 = "Región de Arica y Parinacota"
$ws.Range("P173").Value2 = 300
$ws.Range("Q173").Value2 = 50
$ws.Range("R173").Value2 = "Hortaliza"
